$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sort-sortby")

# --- Capture the existing label text before we move things around ---
$labelText = $ws.Range("H1").Value2

# --- Remove the wrap-text formatting from H1 so it doesn't leak into the
#     new header array formula that will occupy H1:M1 ---
$ws.Range("H1").Style = "Normal"

# --- Clear the old formula blocks (H1:H16 and Q1:V63) ---
$ws.Range("H1:H16").ClearContents()
$ws.Range("Q1:V63").ClearContents()

# --- Block 1 (rows 1-16): headers + SORTBY by emp_last/product now live in H:M ---
$ws.Range("H1:M1").FormulaArray = "=dm_sales_sort[#Headers]"
$ws.Range("H2:M16").FormulaArray = "=SORTBY(dm_sales_sort[], dm_sales_sort[emp_last], -1, dm_sales_sort[product], 1)"

# --- The single-column "trans_id sorted by sales_amt" example moves to column O ---
$ws.Range("O1").Value = $labelText
$ws.Range("O1").WrapText = $true
$ws.Range("O2:O16").FormulaArray = "=SORTBY(dm_sales_sort[trans_id], dm_sales_sort[sales_amt], -1)"

# --- Block 2 (was rows 29-44): SORT by quantity, now rows 19-34 in H:M ---
$ws.Range("H19:M19").FormulaArray = "=dm_sales_sort[#Headers]"
$ws.Range("H20:M34").FormulaArray = "=SORT(dm_sales_sort[], MATCH(""quantity"", dm_sales_sort[#Headers]), -1)"

# --- Block 3 (was rows 48-63): SORTBY by sales_amt, now rows 38-53 in H:M ---
$ws.Range("H38:M38").FormulaArray = "=dm_sales_sort[#Headers]"
$ws.Range("H39:M53").FormulaArray = "=SORTBY(dm_sales_sort[], dm_sales_sort[sales_amt], -1)"

$excel.CalculateFull()

# --- Column widths: old H/I custom widths now belong to O/P; reset H:M to default ---
$ws.Columns.Item(15).ColumnWidth = 19.33
$ws.Columns.Item(16).ColumnWidth = 7.83
$ws.Columns.Item(8).ColumnWidth = 8.43
$ws.Columns.Item(9).ColumnWidth = 8.43
$ws.Columns.Item(10).ColumnWidth = 8.43
$ws.Columns.Item(11).ColumnWidth = 8.43
$ws.Columns.Item(12).ColumnWidth = 8.43
$ws.Columns.Item(13).ColumnWidth = 8.43
